$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 20781
$ws.Range("B2").Value = "Guilherme Nascimento"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45080
$ws.Range("G2").Value = 6426.95

# Row 3
$ws.Range("A3").Value = 69726
$ws.Range("B3").Value = "Laís Casa Grande"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 8853.889999999999

# Row 4
$ws.Range("A4").Value = 92172
$ws.Range("B4").Value = "Asafe Costa"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45083
$ws.Range("G4").Value = 6272.65

# Row 5
$ws.Range("A5").Value = 6339
$ws.Range("B5").Value = "Sr. Rodrigo das Neves"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45099
$ws.Range("G5").Value = 9289.59

# Row 6
$ws.Range("A6").Value = 69213
$ws.Range("B6").Value = "Lucas da Luz"
$ws.Range("C6").Value = "Atendimento ao Cliente"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45089
$ws.Range("G6").Value = 5671.9

# Row 7
$ws.Range("A7").Value = 64630
$ws.Range("B7").Value = "Benjamin Fogaça"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45094
$ws.Range("G7").Value = 8783.780000000001

# Row 8
$ws.Range("A8").Value = 47401
$ws.Range("B8").Value = "Vinicius Borges"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 2381.64

# Row 9
$ws.Range("A9").Value = 60549
$ws.Range("B9").Value = "Sophie Mendes"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Viagem de negocios"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 9417.77

# Row 10
$ws.Range("A10").Value = 30199
$ws.Range("B10").Value = "Maysa Pires"
$ws.Range("C10").Value = "Financeiro"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 45095
$ws.Range("G10").Value = 4954.93

# Row 11
$ws.Range("A11").Value = 4286
$ws.Range("B11").Value = "Sr. Antony da Costa"
$ws.Range("C11").Value = "TI"
$ws.Range("D11").Value = "Consulta medica"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45100
$ws.Range("G11").Value = 5082.97

$wb.Save()
